# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 56 in the "Piña" sheet,
# shifting the existing records (old rows 56-106) down by one row
# (they become rows 57-107).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 56; existing row 56 and below move to 57+.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly record.
$ws.Range("A56").Value = 11
$ws.Range("B56").Value = "Vega Monumental Concepción"
$ws.Range("C56").Value = "Bíobío"
$ws.Range("D56").Value = 44484
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100108
$ws.Range("H56").Value = "Tropicales y subtropicales"
$ws.Range("I56").Value = 100108005
$ws.Range("J56").Value = "Piña"
$ws.Range("K56").Value = "Caramelo"
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 100
$ws.Range("N56").Value = 20000
$ws.Range("O56").Value = 21000
$ws.Range("P56").Value = 20500
$ws.Range("Q56").Value = "`$/caja 12 unidades"
$ws.Range("R56").Value = "Ecuador"
$ws.Range("S56").Value = 1708
$ws.Range("T56").Value = 12
